# Updated cryptos list on Thu Sep  5 15:17:37 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table with
# the latest scrape, and reflects the two rank swaps that happened between
# runs (Polygon/Kaspa and OKB/SuiNetwork traded places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.513.14'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '2.375.64'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").Value = '503.56'
$ws.Range("E5").Value = '  -0.37%  '

$ws.Range("D6").Value = '131.97'
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  -1.41%  '

$ws.Range("D9").Value = '2.381.23'
$ws.Range("E9").Value = '  -2.16%  '

$ws.Range("D10").Value = '0.0985'
$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("E12").Value = '  +2.14%  '

$ws.Range("D13").Value = '4.66'
$ws.Range("E13").Value = '  -0.05%  '

$ws.Range("D14").Value = '2.798.22'
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").Value = '56.451.77'
$ws.Range("E15").Value = '  -1.09%  '

$ws.Range("D16").Value = '21.52'
$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D18").Value = '2.381.39'
$ws.Range("E18").Value = '  -2.09%  '

$ws.Range("D19").Value = '10.06'
$ws.Range("E19").Value = '  -2.21%  '

$ws.Range("D20").Value = '308.87'
$ws.Range("E20").Value = '  -1.24%  '

$ws.Range("D21").Value = '4.02'
$ws.Range("E21").Value = '  -2.05%  '

$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  -4.93%  '

$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").Value = '65.03'
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '0.372'
$ws.Range("E26").Value = '  -2.95%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.150'
$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  -2.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.50'
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").Value = '  -1.72%  '

$ws.Range("E31").Value = '  -2.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  -3.04%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("E34").Value = '  -6.95%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.63%  '

$ws.Range("D36").Value = '17.73'
$ws.Range("E36").Value = '  -1.79%  '

$ws.Range("E37").Value = '  -2.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.80'
$ws.Range("E38").Value = '  -0.80%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.10'
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").Value = '0.804'
$ws.Range("E40").Value = '  -1.79%  '

$ws.Range("D41").Value = '1.42'
$ws.Range("E41").Value = '  -1.99%  '

$ws.Range("D42").Value = '130.42'
$ws.Range("E42").Value = '  -1.22%  '

$ws.Range("E43").Value = '  -0.99%  '

$ws.Range("D44").Value = '4.79'
$ws.Range("E44").Value = '  -4.09%  '

$ws.Range("D45").Value = '0.564'
$ws.Range("E45").Value = '  -0.71%  '

$ws.Range("D46").Value = '0.0908'
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("D47").Value = '245.27'
$ws.Range("E47").Value = '  -4.83%  '

$ws.Range("E48").Value = '  -2.73%  '

$ws.Range("E49").Value = '  -2.01%  '

$ws.Range("D50").Value = '17.03'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").Value = '1.56'
$ws.Range("E51").Value = '  -0.53%  '
